$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.28"
$ws.Range("E2").Value = "'-0.48%"
$ws.Range("D3").Value = "'27.47"
$ws.Range("E3").Value = "'-1.29%"
$ws.Range("D4").Value = "'4.583"
$ws.Range("E4").Value = "'-12.37%"
$ws.Range("E5").Value = "'-1.34%"
$ws.Range("D6").Value = "'6.624"
$ws.Range("E6").Value = "'-1.19%"
$ws.Range("D7").Value = "'0.8580"
$ws.Range("E7").Value = "'-1.44%"
$ws.Range("D8").Value = "'0.9248"
$ws.Range("E8").Value = "'-7.33%"
$ws.Range("D9").Value = "'0.1412"
$ws.Range("E9").Value = "'-1.35%"
$ws.Range("E10").Value = "'-1.32%"
$ws.Range("D11").Value = "'0.07086"
$ws.Range("E11").Value = "'-1.61%"
$ws.Range("D12").Value = "'0.03227"
$ws.Range("E12").Value = "'-0.56%"
$ws.Range("D13").Value = "'0.09183"
$ws.Range("E13").Value = "'-0.67%"
$ws.Range("D14").Value = "'0.001539"
$ws.Range("E14").Value = "'-2.25%"
$ws.Range("D15").Value = "'0.01044"
$ws.Range("E15").Value = "'1,617.70%"
$ws.Range("D16").Value = "'0.006109"
$ws.Range("E16").Value = "'1.96%"
$ws.Range("D17").Value = "'3.518"
$ws.Range("E17").Value = "'0.59%"
$ws.Range("D18").Value = "'3.196"
$ws.Range("E18").Value = "'-2.15%"
$ws.Range("E20").Value = "'-1.32%"
$ws.Range("E21").Value = "'-1.77%"
$ws.Range("D22").Value = "'3.861"
$ws.Range("E22").Value = "'9.42%"
$ws.Range("D23").Value = "'0.04220"
$ws.Range("E23").Value = "'0.96%"
$ws.Range("E24").Value = "'0.04%"
$ws.Range("D25").Value = "'0.004299"
$ws.Range("E25").Value = "'-5.54%"
$ws.Range("E26").Value = "'-0.12%"
$ws.Range("D27").Value = "'0.0001510"
$ws.Range("E27").Value = "'-22.14%"
$ws.Range("D40").Value = "'0.03836"
$ws.Range("E40").Value = "'0.40%"
$ws.Range("D41").Value = "'0.006209"
$ws.Range("E41").Value = "'13.19%"
$ws.Range("D42").Value = "'0.1102"
$ws.Range("E42").Value = "'-0.54%"
$ws.Range("D43").Value = "'0.002200"
$ws.Range("E43").Value = "'-7.71%"
$ws.Range("D44").Value = "'0.01184"
$ws.Range("E44").Value = "'19.09%"
$ws.Range("D45").Value = "'0.00005465"
$ws.Range("E45").Value = "'0.57%"
$ws.Range("D47").Value = "'0.05999"
$ws.Range("E47").Value = "'-45.06%"
$ws.Range("D48").Value = "'0.1215"
$ws.Range("E48").Value = "'5,575.34%"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D50").Value = "'0.0002000"
